# Async Await Working Sample
# Adds a new "WHAT IS A CLASS" section below the existing content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B28").Value = "WHAT IS A CLASS"
$ws.Range("B29").Value = "Class is a type, which we can use to do some work and it also stores some data."
$ws.Range("B31").Value = "In software, we are trying to build the right abstraction and have proper encapsulation and build the right model to manage the complexity"

# Select the last entered cell and scroll the view down, matching the
# author's final view state (topLeftCell A5, active cell B31).
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 5
